# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# --- Sheet "OFF" (first sheet): update row 2 (H) stats ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 203
$wsOff.Range("C2").Value = 149
$wsOff.Range("D2").Value = 49
$wsOff.Range("E2").Value = 14
$wsOff.Range("F2").Value = 4

# --- Sheet "DEF" (second sheet): update row 2 (H) stats ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 155
$wsDef.Range("C2").Value = 108
$wsDef.Range("D2").Value = 35
$wsDef.Range("E2").Value = 14
$wsDef.Range("F2").Value = 3
$wsDef.Range("G2").Value = 1
